# [Abraham]: Fixed name enterprise in file for send email.
# Replace the outdated IFI name "Mutualista Imbabura" with the corrected
# enterprise name in the two data rows of the "Clientes" sheet, wrap the
# text in those cells, and restore the cursor location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newName = "Cooperativa de Ahorro y Crédito San José de Vittoria"

$ws.Range("V2").Value = $newName
$ws.Range("V2").WrapText = $true

$ws.Range("V3").Value = $newName
$ws.Range("V3").WrapText = $true

[void]$ws.Range("Q9").Select()
